$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.902.19'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.464.15'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.92'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.48'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.76'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("E14").Value = '  -2.52%  '
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.589.67'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.463.39'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.71'
$ws.Range("E18").Value = '  -6.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.78'
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.28'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.11'
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.07'
$ws.Range("E24").Value = '  +3.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.88'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '648.48'
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0969'
$ws.Range("E28").Value = '  -2.91%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -2.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.82'
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -3.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.67'
$ws.Range("E36").Value = '  -2.66%  '
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.56'
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.364'
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '149.55'
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("E41").Value = '  -1.76%  '
$ws.Range("E42").Value = '  -2.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0314'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '153.11'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.56'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.24'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.605'
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("E51").Value = '  -1.57%  '
